# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1512.3062  # H40: 1558.5 -> 1512.3062
$ws.Cells.Item(40, 9).Value = 1375.4762  # I40: 1402.9032 -> 1375.4762
$ws.Cells.Item(40, 10).Value = 2333.2856  # J40: 2247.5715 -> 2333.2856
$ws.Cells.Item(40, 11).Value = 1375.4762  # K40: 1402.9032 -> 1375.4762
$ws.Cells.Item(40, 12).Value = 2333.2856  # L40: 2247.5715 -> 2333.2856
$ws.Cells.Item(40, 13).Value = -1200.4762  # M40: -1227.9032 -> -1200.4762
$ws.Cells.Item(40, 14).Value = -2683.2856  # N40: -2597.5715 -> -2683.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 5047.609  # H41: 5709.65 -> 5047.609
$ws.Cells.Item(41, 9).Value = 531.2  # I41: 501.33334 -> 531.2
$ws.Cells.Item(41, 10).Value = 8521.77  # J41: 9971 -> 8521.77
$ws.Cells.Item(41, 11).Value = 531.2  # K41: 501.33334 -> 531.2
$ws.Cells.Item(41, 12).Value = 8521.77  # L41: 9971 -> 8521.77
$ws.Cells.Item(41, 13).Value = -91.20000000000005  # M41: -61.33334000000002 -> -91.20000000000005
$ws.Cells.Item(41, 14).Value = -9401.77  # N41: -10851 -> -9401.77

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 2601.5  # H58: 2814.1 -> 2601.5
$ws.Cells.Item(58, 9).Value = 453.75  # I58: 286.25 -> 453.75
$ws.Cells.Item(58, 10).Value = 4749.25  # J58: 4499.3335 -> 4749.25
$ws.Cells.Item(58, 11).Value = 1361.25  # K58: 858.75 -> 1361.25
$ws.Cells.Item(58, 12).Value = 14247.75  # L58: 13498.0005 -> 14247.75
$ws.Cells.Item(58, 13).Value = -1211.25  # M58: -708.75 -> -1211.25
$ws.Cells.Item(58, 14).Value = -14547.75  # N58: -13798.0005 -> -14547.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 476.81818  # H107: 406.5 -> 476.81818
$ws.Cells.Item(107, 9).Value = 371  # I107: 406.5 -> 371
$ws.Cells.Item(107, 10).Value = 953  # J107: 0 -> 953
$ws.Cells.Item(107, 11).Value = 371  # K107: 406.5 -> 371
$ws.Cells.Item(107, 12).Value = 953  # L107: 0 -> 953
$ws.Cells.Item(107, 13).Value = 1549  # M107: 1513.5 -> 1549
$ws.Cells.Item(107, 14).Value = -4793  # N107: None -> -4793

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 1438.2858  # H112: 1546.5143 -> 1438.2858
$ws.Cells.Item(112, 10).Value = 1450.303  # J112: 1565.091 -> 1450.303
$ws.Cells.Item(112, 12).Value = 4350.909000000001  # L112: 4695.272999999999 -> 4350.909000000001
$ws.Cells.Item(112, 14).Value = -6566.909000000001  # N112: -6911.272999999999 -> -6566.909000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 3537.32  # H113: 3748.739 -> 3537.32
$ws.Cells.Item(113, 9).Value = 2504.3  # I113: 2703.1 -> 2504.3
$ws.Cells.Item(113, 10).Value = 4226  # J113: 4553.077 -> 4226
$ws.Cells.Item(113, 11).Value = 2504.3  # K113: 2703.1 -> 2504.3
$ws.Cells.Item(113, 12).Value = 4226  # L113: 4553.077 -> 4226
$ws.Cells.Item(113, 13).Value = 749.6999999999998  # M113: 550.9000000000001 -> 749.6999999999998
$ws.Cells.Item(113, 14).Value = -10734  # N113: -11061.077 -> -10734

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 100000890  # H125: 45454950 -> 100000890
$ws.Cells.Item(125, 9).Value = 0  # I125: 302.625 -> 0
$ws.Cells.Item(125, 10).Value = 100000890  # J125: 166667340 -> 100000890
$ws.Cells.Item(125, 11).Value = 0  # K125: 2723.625 -> 0
$ws.Cells.Item(125, 12).Value = 900008010  # L125: 1500006060 -> 900008010
$ws.Cells.Item(125, 13).ClearContents()  # M125: -263.625 -> (removed)
$ws.Cells.Item(125, 14).Value = -900012930  # N125: -1500010980 -> -900012930

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 1080.7273  # H129: 955.9286 -> 1080.7273
$ws.Cells.Item(129, 9).Value = 425  # I129: 328.14285 -> 425
$ws.Cells.Item(129, 10).Value = 1455.4286  # J129: 1583.7142 -> 1455.4286
$ws.Cells.Item(129, 11).Value = 1275  # K129: 984.4285500000001 -> 1275
$ws.Cells.Item(129, 12).Value = 4366.2858  # L129: 4751.142599999999 -> 4366.2858
$ws.Cells.Item(129, 13).Value = 3725  # M129: 4015.57145 -> 3725
$ws.Cells.Item(129, 14).Value = -14366.2858  # N129: -14751.1426 -> -14366.2858

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 4230.6816  # H132: 4766.2163 -> 4230.6816
$ws.Cells.Item(132, 9).Value = 2318.0645  # I132: 2678.4 -> 2318.0645
$ws.Cells.Item(132, 10).Value = 8791.538  # J132: 9115.833000000001 -> 8791.538
$ws.Cells.Item(132, 11).Value = 6954.193499999999  # K132: 8035.200000000001 -> 6954.193499999999
$ws.Cells.Item(132, 12).Value = 26374.614  # L132: 27347.499 -> 26374.614
$ws.Cells.Item(132, 13).Value = -4424.193499999999  # M132: -5505.200000000001 -> -4424.193499999999
$ws.Cells.Item(132, 14).Value = -31434.614  # N132: -32407.499 -> -31434.614

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(133, 8).Value = 26593  # H133: 57000 -> 26593
$ws.Cells.Item(133, 10).Value = 26593  # J133: 57000 -> 26593
$ws.Cells.Item(133, 12).Value = 26593  # L133: 57000 -> 26593
$ws.Cells.Item(133, 14).Value = -36713  # N133: -67120 -> -36713

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 862.44446  # H135: 955.1667 -> 862.44446
$ws.Cells.Item(135, 9).Value = 317.30768  # I135: 353 -> 317.30768
$ws.Cells.Item(135, 10).Value = 2279.8  # J135: 2520.8 -> 2279.8
$ws.Cells.Item(135, 11).Value = 2855.76912  # K135: 3177 -> 2855.76912
$ws.Cells.Item(135, 12).Value = 20518.2  # L135: 22687.2 -> 20518.2
$ws.Cells.Item(135, 13).Value = -320.7691199999999  # M135: -642 -> -320.7691199999999
$ws.Cells.Item(135, 14).Value = -25588.2  # N135: -27757.2 -> -25588.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 1204.2778  # H110: 1525.1818 -> 1204.2778
$ws.Cells.Item(110, 9).Value = 1033  # I110: 1242.4286 -> 1033
$ws.Cells.Item(110, 10).Value = 1375.5555  # J110: 2020 -> 1375.5555
$ws.Cells.Item(110, 11).Value = 1033  # K110: 1242.4286 -> 1033
$ws.Cells.Item(110, 12).Value = 1375.5555  # L110: 2020 -> 1375.5555
$ws.Cells.Item(110, 13).Value = 1012  # M110: 802.5714 -> 1012
$ws.Cells.Item(110, 14).Value = -5465.5555  # N110: -6110 -> -5465.5555

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2690.6086  # H122: 2031.4865 -> 2690.6086
$ws.Cells.Item(122, 9).Value = 3062.6667  # I122: 2391.318 -> 3062.6667
$ws.Cells.Item(122, 10).Value = 1993  # J122: 1503.7333 -> 1993
$ws.Cells.Item(122, 11).Value = 9188.000100000001  # K122: 7173.954000000001 -> 9188.000100000001
$ws.Cells.Item(122, 12).Value = 5979  # L122: 4511.199900000001 -> 5979
$ws.Cells.Item(122, 13).Value = -6738.000100000001  # M122: -4723.954000000001 -> -6738.000100000001
$ws.Cells.Item(122, 14).Value = -10879  # N122: -9411.1999 -> -10879

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 1139.4  # H80: 1503.5 -> 1139.4
$ws.Cells.Item(80, 9).Value = 890  # I80: 0 -> 890
$ws.Cells.Item(80, 10).Value = 1201.75  # J80: 1503.5 -> 1201.75
$ws.Cells.Item(80, 11).Value = 890  # K80: 0 -> 890
$ws.Cells.Item(80, 12).Value = 1201.75  # L80: 1503.5 -> 1201.75
$ws.Cells.Item(80, 13).Value = 108  # M80: None -> 108
$ws.Cells.Item(80, 14).Value = -3197.75  # N80: -3499.5 -> -3197.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value = 1139.4  # H83: 1503.5 -> 1139.4
$ws.Cells.Item(83, 9).Value = 890  # I83: 0 -> 890
$ws.Cells.Item(83, 10).Value = 1201.75  # J83: 1503.5 -> 1201.75
$ws.Cells.Item(83, 11).Value = 4450  # K83: 0 -> 4450
$ws.Cells.Item(83, 12).Value = 6008.75  # L83: 7517.5 -> 6008.75
$ws.Cells.Item(83, 13).Value = 542  # M83: None -> 542
$ws.Cells.Item(83, 14).Value = -15992.75  # N83: -17501.5 -> -15992.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 888.7406999999999  # H94: 913.1667 -> 888.7406999999999
$ws.Cells.Item(94, 10).Value = 1026.6666  # J94: 1110 -> 1026.6666
$ws.Cells.Item(94, 12).Value = 1026.6666  # L94: 1110 -> 1026.6666
$ws.Cells.Item(94, 14).Value = -1928.6666  # N94: -2012 -> -1928.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2127.7778  # H107: 2084.842 -> 2127.7778
$ws.Cells.Item(107, 9).Value = 2236.3635  # I107: 2243.5 -> 2236.3635
$ws.Cells.Item(107, 10).Value = 1957.1428  # J107: 1812.8572 -> 1957.1428
$ws.Cells.Item(107, 11).Value = 2236.3635  # K107: 2243.5 -> 2236.3635
$ws.Cells.Item(107, 12).Value = 1957.1428  # L107: 1812.8572 -> 1957.1428
$ws.Cells.Item(107, 13).Value = -316.3634999999999  # M107: -323.5 -> -316.3634999999999
$ws.Cells.Item(107, 14).Value = -5797.1428  # N107: -5652.8572 -> -5797.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 28811.725  # H134: 37904.5 -> 28811.725
$ws.Cells.Item(134, 9).Value = 47107.684  # I134: 68549 -> 47107.684
$ws.Cells.Item(134, 10).Value = 6450  # J134: 7260 -> 6450
$ws.Cells.Item(134, 11).Value = 141323.052  # K134: 205647 -> 141323.052
$ws.Cells.Item(134, 12).Value = 19350  # L134: 21780 -> 19350
$ws.Cells.Item(134, 13).Value = -138788.052  # M134: -203112 -> -138788.052
$ws.Cells.Item(134, 14).Value = -24420  # N134: -26850 -> -24420

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2620.0635  # H31: 3152.0222 -> 2620.0635
$ws.Cells.Item(31, 9).Value = 1612.1052  # I31: 1989.8889 -> 1612.1052
$ws.Cells.Item(31, 10).Value = 4152.16  # J31: 4895.222 -> 4152.16
$ws.Cells.Item(31, 11).Value = 1612.1052  # K31: 1989.8889 -> 1612.1052
$ws.Cells.Item(31, 12).Value = 4152.16  # L31: 4895.222 -> 4152.16
$ws.Cells.Item(31, 13).Value = -1317.1052  # M31: -1694.8889 -> -1317.1052
$ws.Cells.Item(31, 14).Value = -4742.16  # N31: -5485.222 -> -4742.16

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(32, 8).Value = 3004.5  # H32: 4320 -> 3004.5
$ws.Cells.Item(32, 9).Value = 3004.5  # I32: 3184 -> 3004.5
$ws.Cells.Item(32, 10).Value = 0  # J32: 10000 -> 0
$ws.Cells.Item(32, 11).Value = 3004.5  # K32: 3184 -> 3004.5
$ws.Cells.Item(32, 12).Value = 0  # L32: 10000 -> 0
$ws.Cells.Item(32, 13).Value = -2688.5  # M32: -2868 -> -2688.5
$ws.Cells.Item(32, 14).ClearContents()  # N32: -10632 -> (removed)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2620.0635  # H34: 3152.0222 -> 2620.0635
$ws.Cells.Item(34, 9).Value = 1612.1052  # I34: 1989.8889 -> 1612.1052
$ws.Cells.Item(34, 10).Value = 4152.16  # J34: 4895.222 -> 4152.16
$ws.Cells.Item(34, 11).Value = 1612.1052  # K34: 1989.8889 -> 1612.1052
$ws.Cells.Item(34, 12).Value = 4152.16  # L34: 4895.222 -> 4152.16
$ws.Cells.Item(34, 13).Value = -1410.1052  # M34: -1787.8889 -> -1410.1052
$ws.Cells.Item(34, 14).Value = -4556.16  # N34: -5299.222 -> -4556.16

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 29200.865  # H99: 30809.885 -> 29200.865
$ws.Cells.Item(99, 9).Value = 41721.24  # I99: 43426.293 -> 41721.24
$ws.Cells.Item(99, 10).Value = 3116.75  # J99: 3283.182 -> 3116.75
$ws.Cells.Item(99, 11).Value = 41721.24  # K99: 43426.293 -> 41721.24
$ws.Cells.Item(99, 12).Value = 3116.75  # L99: 3283.182 -> 3116.75
$ws.Cells.Item(99, 13).Value = -40223.24  # M99: -41928.293 -> -40223.24
$ws.Cells.Item(99, 14).Value = -6112.75  # N99: -6279.182 -> -6112.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 694.7059  # H107: 657.7143 -> 694.7059
$ws.Cells.Item(107, 9).Value = 311  # I107: 240.8 -> 311
$ws.Cells.Item(107, 10).Value = 1242.8572  # J107: 1700 -> 1242.8572
$ws.Cells.Item(107, 11).Value = 311  # K107: 240.8 -> 311
$ws.Cells.Item(107, 12).Value = 1242.8572  # L107: 1700 -> 1242.8572
$ws.Cells.Item(107, 13).Value = 1609  # M107: 1679.2 -> 1609
$ws.Cells.Item(107, 14).Value = -5082.8572  # N107: -5540 -> -5082.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 29200.865  # H126: 30809.885 -> 29200.865
$ws.Cells.Item(126, 9).Value = 41721.24  # I126: 43426.293 -> 41721.24
$ws.Cells.Item(126, 10).Value = 3116.75  # J126: 3283.182 -> 3116.75
$ws.Cells.Item(126, 11).Value = 125163.72  # K126: 130278.879 -> 125163.72
$ws.Cells.Item(126, 12).Value = 9350.25  # L126: 9849.545999999998 -> 9350.25
$ws.Cells.Item(126, 13).Value = -122693.72  # M126: -127808.879 -> -122693.72
$ws.Cells.Item(126, 14).Value = -14290.25  # N126: -14789.546 -> -14290.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(13, 8).Value = 65  # H13: 76.666664 -> 65
$ws.Cells.Item(13, 9).Value = 65  # I13: 76.666664 -> 65
$ws.Cells.Item(13, 11).Value = 195  # K13: 229.999992 -> 195
$ws.Cells.Item(13, 13).Value = -27  # M13: -61.99999199999999 -> -27

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(16, 8).Value = 216.07692  # H16: 418.18182 -> 216.07692
$ws.Cells.Item(16, 9).Value = 216.07692  # I16: 333.33334 -> 216.07692
$ws.Cells.Item(16, 10).Value = 0  # J16: 800 -> 0
$ws.Cells.Item(16, 11).Value = 648.23076  # K16: 1000.00002 -> 648.23076
$ws.Cells.Item(16, 12).Value = 0  # L16: 2400 -> 0
$ws.Cells.Item(16, 13).Value = -475.23076  # M16: -827.0000200000001 -> -475.23076
$ws.Cells.Item(16, 14).ClearContents()  # N16: -2746 -> (removed)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 2828.125  # H56: 2470.3408 -> 2828.125
$ws.Cells.Item(56, 9).Value = 2828.125  # I56: 2470.3408 -> 2828.125
$ws.Cells.Item(56, 11).Value = 2828.125  # K56: 2470.3408 -> 2828.125
$ws.Cells.Item(56, 13).Value = -2298.125  # M56: -1940.3408 -> -2298.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 150  # H3: 1566.6666 -> 150
$ws.Cells.Item(3, 9).Value = 100  # I3: 300 -> 100
$ws.Cells.Item(3, 10).Value = 200  # J3: 2200 -> 200
$ws.Cells.Item(3, 11).Value = 100  # K3: 300 -> 100
$ws.Cells.Item(3, 12).Value = 200  # L3: 2200 -> 200
$ws.Cells.Item(3, 13).Value = 16  # M3: -184 -> 16
$ws.Cells.Item(3, 14).Value = -432  # N3: -2432 -> -432

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5033.3335  # H70: 3394.8 -> 5033.3335
$ws.Cells.Item(70, 9).Value = 4040  # I70: 3135.7144 -> 4040
$ws.Cells.Item(70, 10).Value = 10000  # J70: 3999.3333 -> 10000
$ws.Cells.Item(70, 11).Value = 4040  # K70: 3135.7144 -> 4040
$ws.Cells.Item(70, 12).Value = 10000  # L70: 3999.3333 -> 10000
$ws.Cells.Item(70, 13).Value = -3770  # M70: -2865.7144 -> -3770
$ws.Cells.Item(70, 14).Value = -10540  # N70: -4539.3333 -> -10540

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 5033.3335  # H73: 3394.8 -> 5033.3335
$ws.Cells.Item(73, 9).Value = 4040  # I73: 3135.7144 -> 4040
$ws.Cells.Item(73, 10).Value = 10000  # J73: 3999.3333 -> 10000
$ws.Cells.Item(73, 11).Value = 4040  # K73: 3135.7144 -> 4040
$ws.Cells.Item(73, 12).Value = 10000  # L73: 3999.3333 -> 10000
$ws.Cells.Item(73, 13).Value = -3104  # M73: -2199.7144 -> -3104
$ws.Cells.Item(73, 14).Value = -11872  # N73: -5871.3333 -> -11872

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1129.5217  # H97: 1327.55 -> 1129.5217
$ws.Cells.Item(97, 9).Value = 798.625  # I97: 904.1667 -> 798.625
$ws.Cells.Item(97, 10).Value = 1885.8572  # J97: 1962.625 -> 1885.8572
$ws.Cells.Item(97, 11).Value = 798.625  # K97: 904.1667 -> 798.625
$ws.Cells.Item(97, 12).Value = 1885.8572  # L97: 1962.625 -> 1885.8572
$ws.Cells.Item(97, 13).Value = -302.625  # M97: -408.1667 -> -302.625
$ws.Cells.Item(97, 14).Value = -2877.8572  # N97: -2954.625 -> -2877.8572

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1887.6666  # H102: 1870.6471 -> 1887.6666
$ws.Cells.Item(102, 9).Value = 2133.0417  # I102: 2141.75 -> 2133.0417
$ws.Cells.Item(102, 10).Value = 1233.3334  # J102: 1220 -> 1233.3334
$ws.Cells.Item(102, 11).Value = 2133.0417  # K102: 2141.75 -> 2133.0417
$ws.Cells.Item(102, 12).Value = 1233.3334  # L102: 1220 -> 1233.3334
$ws.Cells.Item(102, 13).Value = -511.0417000000002  # M102: -519.75 -> -511.0417000000002
$ws.Cells.Item(102, 14).Value = -4477.3334  # N102: -4464 -> -4477.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(110, 8).Value = 38601.5  # H110: 39899.668 -> 38601.5
$ws.Cells.Item(110, 10).Value = 38601.5  # J110: 39899.668 -> 38601.5
$ws.Cells.Item(110, 12).Value = 38601.5  # L110: 39899.668 -> 38601.5
$ws.Cells.Item(110, 14).Value = -46781.5  # N110: -48079.668 -> -46781.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(116, 8).Value = 41000  # H116: 40000 -> 41000
$ws.Cells.Item(116, 10).Value = 41000  # J116: 40000 -> 41000
$ws.Cells.Item(116, 12).Value = 41000  # L116: 40000 -> 41000
$ws.Cells.Item(116, 14).Value = -50178  # N116: -49178 -> -50178

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2607.2593  # H126: 2573.6296 -> 2607.2593
$ws.Cells.Item(126, 9).Value = 2586.0667  # I126: 2411.0557 -> 2586.0667
$ws.Cells.Item(126, 10).Value = 2633.75  # J126: 2898.7778 -> 2633.75
$ws.Cells.Item(126, 11).Value = 7758.2001  # K126: 7233.1671 -> 7758.2001
$ws.Cells.Item(126, 12).Value = 7901.25  # L126: 8696.3334 -> 7901.25
$ws.Cells.Item(126, 13).Value = -5288.2001  # M126: -4763.1671 -> -5288.2001
$ws.Cells.Item(126, 14).Value = -12841.25  # N126: -13636.3334 -> -12841.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(135, 8).Value = 47592.223  # H135: 31965 -> 47592.223
$ws.Cells.Item(135, 10).Value = 47388.332  # J135: 26620 -> 47388.332
$ws.Cells.Item(135, 12).Value = 47388.332  # L135: 26620 -> 47388.332
$ws.Cells.Item(135, 14).Value = -57528.332  # N135: -36760 -> -57528.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 4630  # H32: 3167.5 -> 4630
$ws.Cells.Item(32, 9).Value = 1593.3334  # I32: 890 -> 1593.3334
$ws.Cells.Item(32, 10).Value = 7666.6665  # J32: 10000 -> 7666.6665
$ws.Cells.Item(32, 11).Value = 1593.3334  # K32: 890 -> 1593.3334
$ws.Cells.Item(32, 12).Value = 7666.6665  # L32: 10000 -> 7666.6665
$ws.Cells.Item(32, 13).Value = -1276.3334  # M32: -573 -> -1276.3334
$ws.Cells.Item(32, 14).Value = -8300.666499999999  # N32: -10634 -> -8300.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2872.2856  # H40: 2783.889 -> 2872.2856
$ws.Cells.Item(40, 9).Value = 2711.5  # I40: 2632.5 -> 2711.5
$ws.Cells.Item(40, 11).Value = 2711.5  # K40: 2632.5 -> 2711.5
$ws.Cells.Item(40, 13).Value = -2575.5  # M40: -2496.5 -> -2575.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(33, 8).Value = 4231.5713  # H33: 3716.6667 -> 4231.5713
$ws.Cells.Item(33, 9).Value = 0  # I33: 3500 -> 0
$ws.Cells.Item(33, 10).Value = 4231.5713  # J33: 3760 -> 4231.5713
$ws.Cells.Item(33, 11).Value = 0  # K33: 3500 -> 0
$ws.Cells.Item(33, 12).Value = 4231.5713  # L33: 3760 -> 4231.5713
$ws.Cells.Item(33, 13).ClearContents()  # M33: -3250 -> (removed)
$ws.Cells.Item(33, 14).Value = -4731.5713  # N33: -4260 -> -4731.5713

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(36, 8).Value = 4231.5713  # H36: 3716.6667 -> 4231.5713
$ws.Cells.Item(36, 9).Value = 0  # I36: 3500 -> 0
$ws.Cells.Item(36, 10).Value = 4231.5713  # J36: 3760 -> 4231.5713
$ws.Cells.Item(36, 11).Value = 0  # K36: 3500 -> 0
$ws.Cells.Item(36, 12).Value = 4231.5713  # L36: 3760 -> 4231.5713
$ws.Cells.Item(36, 13).ClearContents()  # M36: -3250 -> (removed)
$ws.Cells.Item(36, 14).Value = -4731.5713  # N36: -4260 -> -4731.5713

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2688.7646  # H96: 2748.3684 -> 2688.7646
$ws.Cells.Item(96, 9).Value = 1268.9  # I96: 1566.5834 -> 1268.9
$ws.Cells.Item(96, 10).Value = 4717.143  # J96: 4774.2856 -> 4717.143
$ws.Cells.Item(96, 11).Value = 1268.9  # K96: 1566.5834 -> 1268.9
$ws.Cells.Item(96, 12).Value = 4717.143  # L96: 4774.2856 -> 4717.143
$ws.Cells.Item(96, 13).Value = 104.0999999999999  # M96: -193.5834 -> 104.0999999999999
$ws.Cells.Item(96, 14).Value = -7463.143  # N96: -7520.2856 -> -7463.143

Write-Host "Edit complete"